# Update cryptos price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.937.97'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '3.147.24'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('D5').Value = "'579.14"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').Value = "'148.54"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.146.99'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = "'37.11"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D15').Value = '3.663.27'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = '64.916.18'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.148.20'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = "'7.13"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = "'502.51"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('D21').Value = "'15.08"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').Value = "'0.713"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('E23').Value = '  -2.69%  '
$ws.Range('D24').Value = "'7.73"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').Value = "'84.17"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').Value = "'9.09"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('D30').Value = "'2.79"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.86%  '
$ws.Range('D31').Value = "'27.44"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'6.48"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.55%  '
$ws.Range('D36').Value = "'54.95"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.82%  '
$ws.Range('D37').Value = "'0.0887"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('D38').Value = "'474.98"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').Value = "'0.0412"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('D40').Value = "'2.92"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('D41').Value = "'8.74"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.84%  '
$ws.Range('D42').Value = '3.002.38'
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = "'2.42"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = "'0.281"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('D46').Value = "'28.23"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.00%  '
$ws.Range('D47').Value = '0.0₃0595'
$ws.Range('E47').Value = '  +2.65%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').Value = "'118.93"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.17%  '
